$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 5556255.5
$ws.Cells.Item(51, 9).Value = 11111111
$ws.Cells.Item(51, 10).Value = 1400
$ws.Cells.Item(51, 11).Value = 11111111
$ws.Cells.Item(51, 12).Value = 1400
$ws.Cells.Item(51, 13).Value = -11110627
$ws.Cells.Item(51, 14).Value = -2368
$ws.Cells.Item(69, 8).Value = 3454.9
$ws.Cells.Item(69, 9).Value = 3183
$ws.Cells.Item(69, 10).Value = 3571.4285
$ws.Cells.Item(69, 11).Value = 9549
$ws.Cells.Item(69, 12).Value = 10714.2855
$ws.Cells.Item(69, 13).Value = -8675
$ws.Cells.Item(69, 14).Value = -12462.2855
$ws.Cells.Item(72, 8).Value = 3454.9
$ws.Cells.Item(72, 9).Value = 3183
$ws.Cells.Item(72, 10).Value = 3571.4285
$ws.Cells.Item(72, 11).Value = 28647
$ws.Cells.Item(72, 12).Value = 32142.8565
$ws.Cells.Item(72, 13).Value = -24279
$ws.Cells.Item(72, 14).Value = -40878.8565
$ws.Cells.Item(112, 8).Value = 1782.375
$ws.Cells.Item(112, 9).Value = 766.6667
$ws.Cells.Item(112, 10).Value = 2391.8
$ws.Cells.Item(112, 11).Value = 2300.0001
$ws.Cells.Item(112, 12).Value = 7175.400000000001
$ws.Cells.Item(112, 13).Value = -1192.0001
$ws.Cells.Item(112, 14).Value = -9391.400000000001
$ws.Cells.Item(124, 8).Value = 35775
$ws.Cells.Item(124, 10).Value = 35775
$ws.Cells.Item(124, 12).Value = 35775
$ws.Cells.Item(124, 14).Value = -45595
$ws.Cells.Item(138, 8).Value = 4633.5713
$ws.Cells.Item(138, 9).Value = 1803.4286
$ws.Cells.Item(138, 10).Value = 5341.107
$ws.Cells.Item(138, 11).Value = 5410.2858
$ws.Cells.Item(138, 12).Value = 16023.321
$ws.Cells.Item(138, 13).Value = -270.2857999999997
$ws.Cells.Item(138, 14).Value = -26303.321
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 17053.87
$ws.Cells.Item(32, 9).Value = 12686.533
$ws.Cells.Item(32, 10).Value = 30155.88
$ws.Cells.Item(32, 11).Value = 12686.533
$ws.Cells.Item(32, 12).Value = 30155.88
$ws.Cells.Item(32, 13).Value = -12399.533
$ws.Cells.Item(32, 14).Value = -30729.88
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2351.3333
$ws.Cells.Item(20, 9).Value = 2421.6
$ws.Cells.Item(20, 11).Value = 2421.6
$ws.Cells.Item(20, 13).Value = -2174.6
$ws.Cells.Item(81, 8).Value = 25027.5
$ws.Cells.Item(81, 10).Value = 25027.5
$ws.Cells.Item(81, 12).Value = 25027.5
$ws.Cells.Item(81, 14).Value = -27149.5
$ws.Cells.Item(84, 8).Value = 25027.5
$ws.Cells.Item(84, 10).Value = 25027.5
$ws.Cells.Item(84, 12).Value = 75082.5
$ws.Cells.Item(84, 14).Value = -85690.5
$ws.Cells.Item(86, 8).Value = 1824.8857
$ws.Cells.Item(86, 9).Value = 1515.2858
$ws.Cells.Item(86, 10).Value = 2289.2856
$ws.Cells.Item(86, 11).Value = 1515.2858
$ws.Cells.Item(86, 12).Value = 2289.2856
$ws.Cells.Item(86, 13).Value = -392.2858000000001
$ws.Cells.Item(86, 14).Value = -4535.2856
$ws.Cells.Item(89, 8).Value = 1824.8857
$ws.Cells.Item(89, 9).Value = 1515.2858
$ws.Cells.Item(89, 10).Value = 2289.2856
$ws.Cells.Item(89, 11).Value = 7576.429
$ws.Cells.Item(89, 12).Value = 11446.428
$ws.Cells.Item(89, 13).Value = -1960.429
$ws.Cells.Item(89, 14).Value = -22678.428
$ws.Cells.Item(134, 8).Value = 6946405.5
$ws.Cells.Item(134, 9).Value = 9260752
$ws.Cells.Item(134, 10).Value = 3366.6667
$ws.Cells.Item(134, 11).Value = 27782256
$ws.Cells.Item(134, 12).Value = 10100.0001
$ws.Cells.Item(134, 13).Value = -27779721
$ws.Cells.Item(134, 14).Value = -15170.0001
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 17860724
$ws.Cells.Item(31, 10).Value = 6185.4614
$ws.Cells.Item(31, 12).Value = 6185.4614
$ws.Cells.Item(31, 14).Value = -6775.4614
$ws.Cells.Item(34, 8).Value = 17860724
$ws.Cells.Item(34, 10).Value = 6185.4614
$ws.Cells.Item(34, 12).Value = 6185.4614
$ws.Cells.Item(34, 14).Value = -6589.4614
$ws.Cells.Item(62, 8).Value = 3364.1
$ws.Cells.Item(62, 9).Value = 3171.25
$ws.Cells.Item(62, 10).Value = 3492.6667
$ws.Cells.Item(62, 11).Value = 3171.25
$ws.Cells.Item(62, 12).Value = 3492.6667
$ws.Cells.Item(62, 13).Value = -2547.25
$ws.Cells.Item(62, 14).Value = -4740.6667
$ws.Cells.Item(65, 8).Value = 3364.1
$ws.Cells.Item(65, 9).Value = 3171.25
$ws.Cells.Item(65, 10).Value = 3492.6667
$ws.Cells.Item(65, 11).Value = 15856.25
$ws.Cells.Item(65, 12).Value = 17463.3335
$ws.Cells.Item(65, 13).Value = -12736.25
$ws.Cells.Item(65, 14).Value = -23703.3335
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 309.94116
$ws.Cells.Item(23, 9).Value = 152.25
$ws.Cells.Item(23, 10).Value = 358.46155
$ws.Cells.Item(23, 11).Value = 456.75
$ws.Cells.Item(23, 12).Value = 1075.38465
$ws.Cells.Item(23, 13).Value = -221.75
$ws.Cells.Item(23, 14).Value = -1545.38465
$ws.Cells.Item(34, 8).Value = 1697.2354
$ws.Cells.Item(34, 9).Value = 187.5
$ws.Cells.Item(34, 10).Value = 2161.7693
$ws.Cells.Item(34, 11).Value = 562.5
$ws.Cells.Item(34, 12).Value = 6485.3079
$ws.Cells.Item(34, 13).Value = -478.5
$ws.Cells.Item(34, 14).Value = -6653.3079
$ws.Cells.Item(39, 8).Value = 4242.8125
$ws.Cells.Item(39, 9).Value = 1995
$ws.Cells.Item(39, 10).Value = 4761.5386
$ws.Cells.Item(39, 11).Value = 5985
$ws.Cells.Item(39, 12).Value = 14284.6158
$ws.Cells.Item(39, 13).Value = -5691
$ws.Cells.Item(39, 14).Value = -14872.6158
$ws.Cells.Item(55, 8).Value = 36338.348
$ws.Cells.Item(55, 9).Value = 111768.86
$ws.Cells.Item(55, 10).Value = 3337.5
$ws.Cells.Item(55, 11).Value = 335306.58
$ws.Cells.Item(55, 12).Value = 10012.5
$ws.Cells.Item(55, 13).Value = -335129.58
$ws.Cells.Item(55, 14).Value = -10366.5
$ws.Cells.Item(93, 8).Value = 11500
$ws.Cells.Item(93, 10).Value = 11500
$ws.Cells.Item(93, 12).Value = 34500
$ws.Cells.Item(93, 14).Value = -38244
$ws.Cells.Item(123, 8).Value = 1388.25
$ws.Cells.Item(123, 9).Value = 1253.3334
$ws.Cells.Item(123, 10).Value = 1433.2222
$ws.Cells.Item(123, 11).Value = 3760.0002
$ws.Cells.Item(123, 12).Value = 4299.6666
$ws.Cells.Item(123, 13).Value = -1310.0002
$ws.Cells.Item(123, 14).Value = -9199.6666
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 8).Value = 224500.56
$ws.Cells.Item(14, 9).Value = 502600
$ws.Cells.Item(14, 10).Value = 2021
$ws.Cells.Item(14, 11).Value = 502600
$ws.Cells.Item(14, 12).Value = 2021
$ws.Cells.Item(14, 13).Value = -502432
$ws.Cells.Item(14, 14).Value = -2357
$ws.Cells.Item(51, 8).Value = 30326
$ws.Cells.Item(51, 10).Value = 30326
$ws.Cells.Item(51, 12).Value = 30326
$ws.Cells.Item(51, 14).Value = -31344
$ws.Cells.Item(132, 8).Value = 7251832.5
$ws.Cells.Item(132, 9).Value = 11116053
$ws.Cells.Item(132, 10).Value = 6418.4375
$ws.Cells.Item(132, 11).Value = 33348159
$ws.Cells.Item(132, 12).Value = 19255.3125
$ws.Cells.Item(132, 13).Value = -33345629
$ws.Cells.Item(132, 14).Value = -24315.3125
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 217.23077
$ws.Cells.Item(55, 9).Value = 70.666664
$ws.Cells.Item(55, 10).Value = 261.2
$ws.Cells.Item(55, 11).Value = 70.666664
$ws.Cells.Item(55, 12).Value = 261.2
$ws.Cells.Item(55, 13).Value = 102.333336
$ws.Cells.Item(55, 14).Value = -607.2
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 256988.75
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 256988.75
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 256988.75
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(4, 14).Value = -257214.75
$ws.Cells.Item(11, 8).Value = 802100.6
$ws.Cells.Item(11, 9).Value = 2000502
$ws.Cells.Item(11, 10).Value = 3166.3333
$ws.Cells.Item(11, 11).Value = 2000502
$ws.Cells.Item(11, 12).Value = 3166.3333
$ws.Cells.Item(11, 13).Value = -2000360
$ws.Cells.Item(11, 14).Value = -3450.3333
$ws.Cells.Item(12, 8).Value = 3850
$ws.Cells.Item(12, 9).Value = 300
$ws.Cells.Item(12, 10).Value = 7400
$ws.Cells.Item(12, 11).Value = 300
$ws.Cells.Item(12, 12).Value = 7400
$ws.Cells.Item(12, 14).Value = -7684
$ws.Cells.Item(75, 8).Value = 29459.916
$ws.Cells.Item(75, 10).Value = 29459.916
$ws.Cells.Item(75, 12).Value = 29459.916
$ws.Cells.Item(75, 14).Value = -31331.916
$ws.Cells.Item(78, 8).Value = 29459.916
$ws.Cells.Item(78, 10).Value = 29459.916
$ws.Cells.Item(78, 12).Value = 88379.74800000001
$ws.Cells.Item(78, 14).Value = -97739.74800000001
$ws.Cells.Item(103, 8).Value = 38866.777
$ws.Cells.Item(103, 10).Value = 38866.777
$ws.Cells.Item(103, 12).Value = 38866.777
$ws.Cells.Item(103, 14).Value = -41210.777
$ws.Cells.Item(133, 8).Value = 39857.5
$ws.Cells.Item(133, 10).Value = 39857.5
$ws.Cells.Item(133, 12).Value = 39857.5
$ws.Cells.Item(133, 14).Value = -49977.5
